# Edit script: apply the commit's changes to the Enron P5 report docx.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the two extra Title paragraphs ("Data Analyst Nanodegree" and
#    "Project P5 - Machine Learning"), keeping "Identify Fraud from Enron Email".
# ---------------------------------------------------------------------
$titleStart = $d.Paragraphs(2).Range.Start
$titleEnd = $d.Paragraphs(3).Range.End
$d.Range($titleStart, $titleEnd).Delete()

# ---------------------------------------------------------------------
# 2) Remove the whole "Data - overview and exploration" section (heading,
#    intro paragraph, 7 bullet stats, observations paragraph, 2 bullet
#    observations) - all the way up to (not including) the next Heading1
#    "Outlier identification and removal".
# ---------------------------------------------------------------------
$secStart = $null
$secEndMarker = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($secStart -eq $null -and $t.StartsWith("Data") -and $t.Contains("overview and exploration")) {
        $secStart = $p.Range.Start
    }
    if ($t -eq "Outlier identification and removal`r") {
        $secEndMarker = $p.Range.Start
        break
    }
}
if ($secStart -ne $null -and $secEndMarker -ne $null) {
    $d.Range($secStart, $secEndMarker).Delete()
}

Write-Output "done"
